$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.876.51"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.821.22"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "3.260.20"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "2.847.59"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.921"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").Value = "51.788.80"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  +7.80%  "
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("B30").Value = "VeChain"
$ws.Range("C30").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +26.81%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0845"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "2.098.16"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.989"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.16%  "
